# Capstone time management.xlsx
# "analysis the source code"
#
# Rows 5 and 6 are new log entries inserted before the existing "Extract
# Source code" row (which used to read "do some change" and sat at row 7),
# row 4's end-time is corrected, and "finish" markers are added in column D
# for the rows that are now complete. Column F already carries the shared
# C-B duration formula all the way down to row 26, so it recalculates on
# its own once B/C are filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: correct the end time and mark it finished ---
$ws.Range("C4").Value = 0.53194444444444444
$ws.Range("D4").Value = "finish"

# --- Row 5 (new): create repository for capstone ---
$ws.Range("A5").Value = "create repository for capstone"
$ws.Range("B5").NumberFormat = "h:mm"
$ws.Range("B5").Value = 0.53472222222222221
$ws.Range("C5").NumberFormat = "h:mm AM/PM"
$ws.Range("C5").Value = 0.54166666666666663
$ws.Range("D5").Value = "finish"

# --- Row 6 (new): Get google source code ---
$ws.Range("A6").Value = "Get google source code"
$ws.Range("B6").NumberFormat = "h:mm AM/PM"
$ws.Range("B6").Value = 0.54166666666666663
$ws.Range("C6").NumberFormat = "h:mm AM/PM"
$ws.Range("C6").Value = 0.625

# --- Row 7: rename to Extract Source code and give it its own times ---
$ws.Range("A7").Value = "Extract Source code"
$ws.Range("B7").NumberFormat = "h:mm AM/PM"
$ws.Range("B7").Value = 0.625
$ws.Range("C7").NumberFormat = "h:mm AM/PM"
$ws.Range("C7").Value = 0.875

# Selection moved to D7 as part of this edit.
$ws.Range("D7").Select()
